$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 183
$ws.Range("I2").Value = 556
$ws.Range("J2").Value = 2331
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 582
$ws.Range("M2").Value = 49
$ws.Range("N2").Value = 386
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 244
$ws.Range("T2").Value = 413
$ws.Range("U2").Value = 34
$ws.Range("V2").Value = 3476
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3548
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 21

$wb.Save()
